# ImportarUsuarios.xlsx — add "Unidade Escola" (unidade_escola) list column
# to Aluno / Coordenador / Pais / Professor sheets, backed by a new
# Auxiliar!H2:H4 list + "unidade_escola" defined name.

$wb = $excel.ActiveWorkbook

$wsAluno       = $wb.Worksheets.Item("Aluno")
$wsCoordenador = $wb.Worksheets.Item("Coordenador")
$wsPais        = $wb.Worksheets.Item("Pais")
$wsProfessor   = $wb.Worksheets.Item("Professor")
$wsAuxiliar    = $wb.Worksheets.Item("Auxiliar")

# ---------------------------------------------------------------------
# Auxiliar sheet: new column H holding the "Unidade Escola" list values
# ---------------------------------------------------------------------
$wsAuxiliar.Range("G1").Copy()
$wsAuxiliar.Range("H1").PasteSpecial(-4122)
$wsAuxiliar.Range("H1").Value = "Unidade Escola"
$wsAuxiliar.Range("H2").Value = "Jaguariúna"
$wsAuxiliar.Range("H3").Value = "Santo Antônio de Posse"
$wsAuxiliar.Range("H4").Value = "Ambos"
$wsAuxiliar.Columns.Item(8).ColumnWidth = 18.65

# ---------------------------------------------------------------------
# Workbook-level defined name used by the new data validations
# ---------------------------------------------------------------------
$wb.Names.Add("unidade_escola", "=Auxiliar!`$H`$2:`$H`$4")

# ---------------------------------------------------------------------
# Aluno (sheet1): new column X
# ---------------------------------------------------------------------
$wsAluno.Range("W1").Copy()
$wsAluno.Range("X1").PasteSpecial(-4122)
$wsAluno.Range("X1").Value = "Unidade Escola"
$wsAluno.Range("X2").Value = "Jaguariúna"
$wsAluno.Range("X2").Validation.Add(3, 1, 1, "unidade_escola")
$wsAluno.Columns.Item(24).ColumnWidth = 14.15

# ---------------------------------------------------------------------
# Coordenador (sheet2): new column S
# ---------------------------------------------------------------------
$wsCoordenador.Range("R1").Copy()
$wsCoordenador.Range("S1").PasteSpecial(-4122)
$wsCoordenador.Range("S1").Value = "Unidade Escola"
$wsCoordenador.Range("S2").Value = "Jaguariúna"
$wsCoordenador.Range("S2").Validation.Add(3, 1, 1, "unidade_escola")
$wsCoordenador.Columns.Item(19).ColumnWidth = 15.65

# ---------------------------------------------------------------------
# Pais (sheet3): new column X
# ---------------------------------------------------------------------
$wsPais.Range("W1").Copy()
$wsPais.Range("X1").PasteSpecial(-4122)
$wsPais.Range("X1").Value = "Unidade Escola"
$wsPais.Range("X2").Value = "Ambos"
$wsPais.Range("X2").Validation.Add(3, 1, 1, "unidade_escola")
$wsPais.Columns.Item(24).ColumnWidth = 13.65

# ---------------------------------------------------------------------
# Professor (sheet4): new column S
# ---------------------------------------------------------------------
$wsProfessor.Range("R1").Copy()
$wsProfessor.Range("S1").PasteSpecial(-4122)
$wsProfessor.Range("S1").Value = "Unidade Escola"
$wsProfessor.Range("S2").Value = "Ambos"
$wsProfessor.Range("S2").Validation.Add(3, 1, 1, "unidade_escola")
$wsProfessor.Columns.Item(19).ColumnWidth = 13.65

# ---------------------------------------------------------------------
# Selection bookkeeping
# ---------------------------------------------------------------------
[void]$wsAluno.Range("A1").Select()
[void]$wsCoordenador.Range("A1").Select()
[void]$wsPais.Range("A1").Select()
[void]$wsProfessor.Range("A1").Select()
[void]$wsAuxiliar.Range("A2").Select()

[void]$wsAluno.Select()
[void]$wsAluno.Range("A1").Select()
